$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 1.030866699911229
$ws.Cells.Item(2, 4).Value = 1.046342341077328
$ws.Cells.Item(2, 5).Value = 1.041946786723488
$ws.Cells.Item(2, 6).Value = 1.05023049916031
$ws.Cells.Item(2, 9).Value = 1.061817213719469
$ws.Cells.Item(2, 10).Value = 1.052161820581608
$ws.Cells.Item(2, 11).Value = 1.057191709460576
$ws.Cells.Item(2, 12).Value = 1.052851104093847
$ws.Cells.Item(2, 13).Value = 1.061031738687747
$ws.Cells.Item(2, 14).Value = 1.020744219655744
$ws.Cells.Item(2, 16).Value = 1.056874736064959

$ws.Cells.Item(3, 3).Value = 1.035061613573453
$ws.Cells.Item(3, 4).Value = 1.049417894264705
$ws.Cells.Item(3, 5).Value = 1.045258157214057
$ws.Cells.Item(3, 6).Value = 1.053383853135484
$ws.Cells.Item(3, 9).Value = 1.063219606452499
$ws.Cells.Item(3, 10).Value = 1.054640047940227
$ws.Cells.Item(3, 11).Value = 1.059455223527049
$ws.Cells.Item(3, 12).Value = 1.055342936210032
$ws.Cells.Item(3, 13).Value = 1.063376372964199
$ws.Cells.Item(3, 14).Value = 1.021606556811344
$ws.Cells.Item(3, 16).Value = 1.058730313716159

$ws.Cells.Item(4, 3).Value = 1.037724209853263
$ws.Cells.Item(4, 4).Value = 1.051372535408231
$ws.Cells.Item(4, 5).Value = 1.047365353551287
$ws.Cells.Item(4, 6).Value = 1.05539162556737
$ws.Cells.Item(4, 9).Value = 1.064099714190244
$ws.Cells.Item(4, 10).Value = 1.056209791180171
$ws.Cells.Item(4, 11).Value = 1.060888088564982
$ws.Cells.Item(4, 12).Value = 1.056923850348192
$ws.Cells.Item(4, 13).Value = 1.064864508191894
$ws.Cells.Item(4, 14).Value = 1.022152330839682
$ws.Cells.Item(4, 16).Value = 1.059908043246595

$ws.Cells.Item(5, 3).Value = 1.038831554947938
$ws.Cells.Item(5, 4).Value = 1.052186012473465
$ws.Cells.Item(5, 5).Value = 1.048242972475436
$ws.Cells.Item(5, 6).Value = 1.056228101687807
$ws.Cells.Item(5, 9).Value = 1.064463300156339
$ws.Cells.Item(5, 10).Value = 1.056861821799077
$ws.Cells.Item(5, 11).Value = 1.061483045555445
$ws.Cells.Item(5, 12).Value = 1.057581133787154
$ws.Cells.Item(5, 13).Value = 1.065483362177047
$ws.Cells.Item(5, 14).Value = 1.022378919737394
$ws.Cells.Item(5, 16).Value = 1.060397811792501

$ws.Cells.Item(6, 3).Value = 1.039016789758203
$ws.Cells.Item(6, 4).Value = 1.052322121395655
$ws.Cells.Item(6, 5).Value = 1.048389851998415
$ws.Cells.Item(6, 6).Value = 1.05636811097315
$ws.Cells.Item(6, 9).Value = 1.064523975617917
$ws.Cells.Item(6, 10).Value = 1.056970843927647
$ws.Cells.Item(6, 11).Value = 1.061582511703987
$ws.Cells.Item(6, 12).Value = 1.057691070238008
$ws.Cells.Item(6, 13).Value = 1.065586879407901
$ws.Cells.Item(6, 14).Value = 1.022416799597972
$ws.Cells.Item(6, 16).Value = 1.06047973655363

$ws.Cells.Item(7, 3).Value = 1.037739052941829
$ws.Cells.Item(7, 4).Value = 1.051383437260336
$ws.Cells.Item(7, 5).Value = 1.047377112419384
$ws.Cells.Item(7, 6).Value = 1.055402832136574
$ws.Cells.Item(7, 9).Value = 1.06410459746403
$ws.Cells.Item(7, 10).Value = 1.056218534372049
$ws.Cells.Item(7, 11).Value = 1.060896067318202
$ws.Cells.Item(7, 12).Value = 1.056932661554719
$ws.Cells.Item(7, 13).Value = 1.064872803664326
$ws.Cells.Item(7, 14).Value = 1.022155369657141
$ws.Cells.Item(7, 16).Value = 1.05991460838609

$ws.Cells.Item(8, 3).Value = 1.032295354628655
$ws.Cells.Item(8, 4).Value = 1.04738922124182
$ws.Cells.Item(8, 5).Value = 1.043073379310918
$ws.Cells.Item(8, 6).Value = 1.051303091828249
$ws.Cells.Item(8, 9).Value = 1.062296900360233
$ws.Cells.Item(8, 10).Value = 1.05300648495905
$ws.Cells.Item(8, 11).Value = 1.057963371173225
$ws.Cells.Item(8, 12).Value = 1.053699869641059
$ws.Cells.Item(8, 13).Value = 1.061830237141576
$ws.Cells.Item(8, 14).Value = 1.02103822302802
$ws.Cells.Item(8, 16).Value = 1.057506679853505

$ws.Cells.Item(9, 3).Value = 1.022287803511089
$ws.Cells.Item(9, 4).Value = 1.040068478969739
$ws.Cells.Item(9, 5).Value = 1.035206223159341
$ws.Cells.Item(9, 6).Value = 1.043818011451999
$ws.Cells.Item(9, 9).Value = 1.058895892047574
$ws.Cells.Item(9, 10).Value = 1.047077347208076
$ws.Cells.Item(9, 11).Value = 1.052543281414122
$ws.Cells.Item(9, 12).Value = 1.047752808177457
$ws.Cells.Item(9, 13).Value = 1.056238037922602
$ws.Cells.Item(9, 14).Value = 1.018972804097555
$ws.Cells.Item(9, 16).Value = 1.053080917544924

$ws.Cells.Item(10, 3).Value = 1.015395520541208
$ws.Cells.Item(10, 4).Value = 1.035065420094653
$ws.Cells.Item(10, 5).Value = 1.029848335571778
$ws.Cells.Item(10, 6).Value = 1.038756068642721
$ws.Cells.Item(10, 9).Value = 1.056526926490424
$ws.Cells.Item(10, 10).Value = 1.043011352702047
$ws.Cells.Item(10, 11).Value = 1.048829161223509
$ws.Cells.Item(10, 12).Value = 1.043698800670906
$ws.Cells.Item(10, 13).Value = 1.052459136171366
$ws.Cells.Item(10, 14).Value = 1.017568156621104
$ws.Cells.Item(10, 16).Value = 1.050140935873401

$ws.Cells.Item(11, 3).Value = 1.013026725183208
$ws.Cells.Item(11, 4).Value = 1.033521105437651
$ws.Cells.Item(11, 5).Value = 1.028240255915968
$ws.Cells.Item(11, 6).Value = 1.037471964494867
$ws.Cells.Item(11, 9).Value = 1.055887623373812
$ws.Cells.Item(11, 10).Value = 1.041867442056493
$ws.Cells.Item(11, 11).Value = 1.04783738956066
$ws.Cells.Item(11, 12).Value = 1.042648554741145
$ws.Cells.Item(11, 13).Value = 1.051720199611142
$ws.Cells.Item(11, 14).Value = 1.01728909773634
$ws.Cells.Item(11, 16).Value = 1.04998664981317

$ws.Cells.Item(12, 3).Value = 1.012402064655563
$ws.Cells.Item(12, 4).Value = 1.033197795160395
$ws.Cells.Item(12, 5).Value = 1.027927721119286
$ws.Cells.Item(12, 6).Value = 1.03734860079428
$ws.Cells.Item(12, 9).Value = 1.055807038913233
$ws.Cells.Item(12, 10).Value = 1.041690665760029
$ws.Cells.Item(12, 11).Value = 1.047715466090333
$ws.Cells.Item(12, 12).Value = 1.042538740592518
$ws.Cells.Item(12, 13).Value = 1.051793613212225
$ws.Cells.Item(12, 14).Value = 1.017321746842973
$ws.Cells.Item(12, 16).Value = 1.050367231470598

$ws.Cells.Item(13, 3).Value = 1.013082327719042
$ws.Cells.Item(13, 4).Value = 1.033796582037213
$ws.Cells.Item(13, 5).Value = 1.028594971051661
$ws.Cells.Item(13, 6).Value = 1.038121017341536
$ws.Cells.Item(13, 9).Value = 1.056158129222564
$ws.Cells.Item(13, 10).Value = 1.04225252108183
$ws.Cells.Item(13, 11).Value = 1.048261764751412
$ws.Cells.Item(13, 12).Value = 1.043151836566068
$ws.Cells.Item(13, 13).Value = 1.052510883228845
$ws.Cells.Item(13, 14).Value = 1.017598042818703
$ws.Cells.Item(13, 16).Value = 1.051208094173664

$ws.Cells.Item(14, 3).Value = 1.01414611554473
$ws.Cells.Item(14, 4).Value = 1.034630171939516
$ws.Cells.Item(14, 5).Value = 1.029502099182135
$ws.Cells.Item(14, 6).Value = 1.039064396537775
$ws.Cells.Item(14, 9).Value = 1.056597064550698
$ws.Cells.Item(14, 10).Value = 1.042976846269662
$ws.Cells.Item(14, 11).Value = 1.048943438540026
$ws.Cells.Item(14, 12).Value = 1.043904568783201
$ws.Cells.Item(14, 13).Value = 1.053301446736323
$ws.Cells.Item(14, 14).Value = 1.017899826627088
$ws.Cells.Item(14, 16).Value = 1.052004720816475

$ws.Cells.Item(15, 3).Value = 1.014693715690952
$ws.Cells.Item(15, 4).Value = 1.035039905725823
$ws.Cells.Item(15, 5).Value = 1.029943401670102
$ws.Cells.Item(15, 6).Value = 1.039500581818946
$ws.Cells.Item(15, 9).Value = 1.05680191294151
$ws.Cells.Item(15, 10).Value = 1.043320270118034
$ws.Cells.Item(15, 11).Value = 1.049261519472123
$ws.Cells.Item(15, 12).Value = 1.044252972477862
$ws.Cells.Item(15, 13).Value = 1.053646125428694
$ws.Cells.Item(15, 14).Value = 1.018029087140139
$ws.Cells.Item(15, 16).Value = 1.052314370357411

$ws.Cells.Item(16, 3).Value = 1.017484356664727
$ws.Cells.Item(16, 4).Value = 1.037048121475264
$ws.Cells.Item(16, 5).Value = 1.03208714931664
$ws.Cells.Item(16, 6).Value = 1.041518620496329
$ws.Cells.Item(16, 9).Value = 1.057755017006092
$ws.Cells.Item(16, 10).Value = 1.044946421567343
$ws.Cells.Item(16, 11).Value = 1.050744387101237
$ws.Cells.Item(16, 12).Value = 1.04586521013307
$ws.Cells.Item(16, 13).Value = 1.055142034817269
$ws.Cells.Item(16, 14).Value = 1.018574129867487
$ws.Cells.Item(16, 16).Value = 1.05345838053969

$ws.Cells.Item(17, 3).Value = 1.019067684776379
$ws.Cells.Item(17, 4).Value = 1.038151577746192
$ws.Cells.Item(17, 5).Value = 1.033256514298524
$ws.Cells.Item(17, 6).Value = 1.042568501851048
$ws.Cells.Item(17, 9).Value = 1.058252797412528
$ws.Cells.Item(17, 10).Value = 1.045812338429825
$ws.Cells.Item(17, 11).Value = 1.051522179921314
$ws.Cells.Item(17, 12).Value = 1.04670552381048
$ws.Cells.Item(17, 13).Value = 1.055869164133548
$ws.Cells.Item(17, 14).Value = 1.0188333224065
$ws.Cells.Item(17, 16).Value = 1.053905362931085

$ws.Cells.Item(18, 3).Value = 1.019747972182134
$ws.Cells.Item(18, 4).Value = 1.038560145625778
$ws.Cells.Item(18, 5).Value = 1.033672739698001
$ws.Cells.Item(18, 6).Value = 1.042846041184934
$ws.Cells.Item(18, 9).Value = 1.058394282131542
$ws.Cells.Item(18, 10).Value = 1.046085138436412
$ws.Cells.Item(18, 11).Value = 1.051744896719176
$ws.Cells.Item(18, 12).Value = 1.046934493553546
$ws.Cells.Item(18, 13).Value = 1.055964025355469
$ws.Cells.Item(18, 14).Value = 1.018860557769241
$ws.Cells.Item(18, 16).Value = 1.053745335805444

$ws.Cells.Item(19, 3).Value = 1.019606090809323
$ws.Cells.Item(19, 4).Value = 1.038334913058711
$ws.Cells.Item(19, 5).Value = 1.033401408772226
$ws.Cells.Item(19, 6).Value = 1.042417775851845
$ws.Cells.Item(19, 9).Value = 1.058212430884697
$ws.Cells.Item(19, 10).Value = 1.045817908083989
$ws.Cells.Item(19, 11).Value = 1.05146201557138
$ws.Cells.Item(19, 12).Value = 1.046605874758073
$ws.Cells.Item(19, 13).Value = 1.055481561129864
$ws.Cells.Item(19, 14).Value = 1.018680435138518
$ws.Cells.Item(19, 16).Value = 1.053042495665236

$ws.Cells.Item(20, 2).Value = 1.05
$ws.Cells.Item(20, 3).Value = 1.017185619175808
$ws.Cells.Item(20, 4).Value = 1.036364371199725
$ws.Cells.Item(20, 5).Value = 1.031238678447608
$ws.Cells.Item(20, 6).Value = 1.040070707566276
$ws.Cells.Item(20, 9).Value = 1.057146789938635
$ws.Cells.Item(20, 10).Value = 1.044069768725718
$ws.Cells.Item(20, 11).Value = 1.049796577063583
$ws.Cells.Item(20, 12).Value = 1.044753607748218
$ws.Cells.Item(20, 13).Value = 1.053443747739014
$ws.Cells.Item(20, 14).Value = 1.017934474555378
$ws.Cells.Item(20, 16).Value = 1.050909738534129

$ws.Cells.Item(21, 3).Value = 1.011872751819274
$ws.Cells.Item(21, 4).Value = 1.032482243640313
$ws.Cells.Item(21, 5).Value = 1.02707763026889
$ws.Cells.Item(21, 6).Value = 1.036095895903691
$ws.Cells.Item(21, 9).Value = 1.055268564922577
$ws.Cells.Item(21, 10).Value = 1.040882330547958
$ws.Cells.Item(21, 11).Value = 1.046873297946104
$ws.Cells.Item(21, 12).Value = 1.041563547450602
$ws.Cells.Item(21, 13).Value = 1.050424236776499
$ws.Cells.Item(21, 14).Value = 1.016811570183602
$ws.Cells.Item(21, 16).Value = 1.048479748137088

$ws.Cells.Item(22, 3).Value = 1.008501980331464
$ws.Cells.Item(22, 4).Value = 1.030034537649771
$ws.Cells.Item(22, 5).Value = 1.024459354346946
$ws.Cells.Item(22, 6).Value = 1.03361096608963
$ws.Cells.Item(22, 9).Value = 1.054078603901536
$ws.Cells.Item(22, 10).Value = 1.038873495982491
$ws.Cells.Item(22, 11).Value = 1.04503359601827
$ws.Cells.Item(22, 12).Value = 1.039561355581607
$ws.Cells.Item(22, 13).Value = 1.048544748246003
$ws.Cells.Item(22, 14).Value = 1.016110272612222
$ws.Cells.Item(22, 16).Value = 1.046992273718561

$ws.Cells.Item(23, 3).Value = 1.010295977014325
$ws.Cells.Item(23, 4).Value = 1.031336788190511
$ws.Cells.Item(23, 5).Value = 1.025852104732832
$ws.Cells.Item(23, 6).Value = 1.034932641384527
$ws.Cells.Item(23, 9).Value = 1.054712814895834
$ws.Cells.Item(23, 10).Value = 1.039942834936753
$ws.Cells.Item(23, 11).Value = 1.046012960524186
$ws.Cells.Item(23, 12).Value = 1.040626882746103
$ws.Cells.Item(23, 13).Value = 1.049544901850732
$ws.Cells.Item(23, 14).Value = 1.016483609616823
$ws.Cells.Item(23, 16).Value = 1.047783820768209

$ws.Cells.Item(24, 3).Value = 1.017205245570787
$ws.Cells.Item(24, 4).Value = 1.036362077771769
$ws.Cells.Item(24, 5).Value = 1.031232105892387
$ws.Cells.Item(24, 6).Value = 1.040041183784054
$ws.Cells.Item(24, 9).Value = 1.057135701328023
$ws.Cells.Item(24, 10).Value = 1.044056709573273
$ws.Cells.Item(24, 11).Value = 1.049779369888403
$ws.Cells.Item(24, 12).Value = 1.044732085705134
$ws.Cells.Item(24, 13).Value = 1.053399817494832
$ws.Cells.Item(24, 14).Value = 1.017919327364111
$ws.Cells.Item(24, 16).Value = 1.050834691345419

$ws.Cells.Item(25, 3).Value = 1.024926880576954
$ws.Cells.Item(25, 4).Value = 1.041996138398976
$ws.Cells.Item(25, 5).Value = 1.037275371026638
$ws.Cells.Item(25, 6).Value = 1.045785557973607
$ws.Cells.Item(25, 9).Value = 1.059801606450214
$ws.Cells.Item(25, 10).Value = 1.048643521698025
$ws.Cells.Item(25, 11).Value = 1.053975710994689
$ws.Cells.Item(25, 12).Value = 1.049321335864074
$ws.Cells.Item(25, 13).Value = 1.057712383682184
$ws.Cells.Item(25, 14).Value = 1.01951872847922
$ws.Cells.Item(25, 16).Value = 1.05424774255203
